# Apply the edits described by the diff:
# - A31 value changes from 1 to 2
# - A33 value changes from 28 to 29
# The dependent SUM formulas in A34 and A62 will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A31").Value = 2
$ws.Range("A33").Value = 29

$excel.CalculateFullRebuild()
$wb.Save()
